$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 127, pushing existing rows 127-130 down to 128-131
$ws.Rows("127").Insert()

# Fill in the new row 127 with the new weekly data point
$ws.Range("A127").Value = 5
$ws.Range("B127").Value = "Macroferia Regional de Talca"
$ws.Range("C127").Value = "Maule"
$ws.Range("D127").Value = 45239
$ws.Range("E127").Value = 7
$ws.Range("F127").Value = 300000000
$ws.Range("G127").Value = "Espárragos"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 4000
$ws.Range("K127").Value = 1400
$ws.Range("L127").Value = 1400
$ws.Range("M127").Value = 1400
$ws.Range("N127").Value = "$/kilo"
$ws.Range("O127").Value = "Provincia de Linares"
$ws.Range("P127").Value = 1400
$ws.Range("Q127").Value = 1
$ws.Range("R127").Value = "Hortaliza"
